$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the now-redundant "duplicate default" cell style (old cellXfs index 1) ---
# Reset every cell/column/row that used to carry that duplicate style back to the
# workbook's true "Normal" style so the style table can collapse 5 -> 4 entries.
$ws.Columns("B:C").Style = "Normal"
$ws.Range("B1").Style = "Normal"
$ws.Range("C1").Style = "Normal"
$ws.Range("H1").Style = "Normal"
$ws.Range("H2").Style = "Normal"
$ws.Range("A5").Style = "Normal"
$ws.Range("A6").Style = "Normal"
$ws.Range("A8").Style = "Normal"
$ws.Range("A11").Style = "Normal"
$ws.Range("D15:G15,I15").Style = "Normal"
$ws.Range("D16:G16,I16").Style = "Normal"
$ws.Range("D17:G17,I17").Style = "Normal"

# Row 7 carried the duplicate style at the row level (customFormat) as well as on A7.
$ws.Rows("7:7").ClearFormats()

# Rows 13 and 14 held no real data, only the duplicate row-level style - remove them
# outright (delete + reinsert keeps every other row number stable).
$ws.Rows("13:14").Delete()
$ws.Rows("13:14").Insert()

# --- Add the new SVR parameter columns (K, L, M) loaded from pred_par ---
$ws.Range("K1").Value = "svr_kernel_scale"
$ws.Range("L1").Value = "svr_epsilon"
$ws.Range("M1").Value = "svr_box_constraint"
$ws.Range("K2").Value = 100
$ws.Range("L2").Value = 0.2
$ws.Range("M2").Value = 20

# Matches the saved selection left behind in the authored workbook.
$ws.Range("J10").Select()
